# Generate Report for Archive
#
# 1. Update the "Ready for handoff" status text to "In Translation" everywhere
#    it appears (Overview sheet columns E/F, and the "Status" column (C) on
#    the zh-cn and de-de language sheets).
# 2. Shrink the "Status" column width on those same columns to match the
#    new (shorter) text, matching the narrower column width used in the
#    updated report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status text -------------------------------------------------

if ($wsOverview.Range("E2").Value2 -eq $oldStatus) { $wsOverview.Range("E2").Value = $newStatus }
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) { $wsOverview.Range("F2").Value = $newStatus }
if ($wsZhCn.Range("C2").Value2 -eq $oldStatus) { $wsZhCn.Range("C2").Value = $newStatus }
if ($wsDeDe.Range("C2").Value2 -eq $oldStatus) { $wsDeDe.Range("C2").Value = $newStatus }

# --- Resize the columns that held the status text ---------------------------
# The report now shows the narrower "In Translation" text instead of
# "Ready for handoff", so the previously auto-fitted Status columns need to
# shrink to match (from ~17.22 chars wide down to ~13.41 chars wide).
#
# Excel's ColumnWidth is expressed in "characters" but gets stored (and
# re-read) after being snapped to the worksheet's underlying pixel grid, so
# we solve for the input that lands on the closest achievable stored width.

$desiredStoredWidth = 13.4101845877511
$columnWidthInput = ($desiredStoredWidth * 6 - 5) / 6

$wsOverview.Columns.Item(5).ColumnWidth = $columnWidthInput
$wsOverview.Columns.Item(6).ColumnWidth = $columnWidthInput
$wsZhCn.Columns.Item(3).ColumnWidth = $columnWidthInput
$wsDeDe.Columns.Item(3).ColumnWidth = $columnWidthInput
